$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 243.95833
$ws.Range("I9").Value = 393.2857
$ws.Range("J9").Value = 182.47058
$ws.Range("K9").Value = 393.2857
$ws.Range("L9").Value = 182.47058
$ws.Range("M9").Value = -224.2857
$ws.Range("N9").Value = -520.47058
$ws.Range("H113").Value = 6319.7334
$ws.Range("I113").Value = 7350.125
$ws.Range("J113").Value = 5142.143
$ws.Range("K113").Value = 7350.125
$ws.Range("L113").Value = 5142.143
$ws.Range("M113").Value = -4096.125
$ws.Range("N113").Value = -11650.143
$ws.Range("H115").Value = 1406.75
$ws.Range("I115").Value = 1406.75
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 4220.25
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -2653.25
$ws.Range("N115").ClearContents()
$ws.Range("H132").Value = 1966.3077
$ws.Range("I132").Value = 1519.0555
$ws.Range("K132").Value = 4557.166499999999
$ws.Range("M132").Value = -2027.166499999999
$ws.Range("H137").Value = 2962.1
$ws.Range("I137").Value = 2944.923
$ws.Range("J137").Value = 2994
$ws.Range("K137").Value = 8834.769
$ws.Range("L137").Value = 8982
$ws.Range("M137").Value = -6284.769
$ws.Range("N137").Value = -14082
$ws.Range("H138").Value = 2954
$ws.Range("I138").Value = 1865.25
$ws.Range("J138").Value = 3722.5293
$ws.Range("K138").Value = 5595.75
$ws.Range("L138").Value = 11167.5879
$ws.Range("M138").Value = -455.75
$ws.Range("N138").Value = -21447.5879
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 30270.762
$ws.Range("I32").Value = 6091
$ws.Range("K32").Value = 6091
$ws.Range("M32").Value = -5804
$ws.Range("H38").Value = 3277.25
$ws.Range("I38").Value = 3736.6667
$ws.Range("J38").Value = 1899
$ws.Range("K38").Value = 3736.6667
$ws.Range("L38").Value = 1899
$ws.Range("M38").Value = -3269.6667
$ws.Range("N38").Value = -2833
$ws.Range("H61").Value = 1907.6863
$ws.Range("I61").Value = 1780.8837
$ws.Range("J61").Value = 2589.25
$ws.Range("K61").Value = 1780.8837
$ws.Range("L61").Value = 2589.25
$ws.Range("M61").Value = -1568.8837
$ws.Range("N61").Value = -3013.25
$ws.Range("H110").Value = 1978.5
$ws.Range("I110").Value = 1865
$ws.Range("K110").Value = 1865
$ws.Range("M110").Value = 180
$ws.Range("H132").Value = 6441.7437
$ws.Range("I132").Value = 6711.5947
$ws.Range("J132").Value = 1449.5
$ws.Range("K132").Value = 20134.7841
$ws.Range("L132").Value = 4348.5
$ws.Range("M132").Value = -17604.7841
$ws.Range("N132").Value = -9408.5
$ws.Range("H136").Value = 1907.6863
$ws.Range("I136").Value = 1780.8837
$ws.Range("J136").Value = 2589.25
$ws.Range("K136").Value = 5342.6511
$ws.Range("L136").Value = 7767.75
$ws.Range("M136").Value = -2792.6511
$ws.Range("N136").Value = -12867.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 25000
$ws.Range("J9").Value = 25000
$ws.Range("L9").Value = 25000
$ws.Range("N9").Value = -25336
$ws.Range("H134").Value = 1886.697
$ws.Range("I134").Value = 1718.72
$ws.Range("J134").Value = 2411.625
$ws.Range("K134").Value = 5156.16
$ws.Range("L134").Value = 7234.875
$ws.Range("M134").Value = -2621.16
$ws.Range("N134").Value = -12304.875
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 25641322
$ws.Range("J7").Value = 594.61536
$ws.Range("L7").Value = 594.61536
$ws.Range("N7").Value = -820.61536
$ws.Range("H16").Value = 1425.1428
$ws.Range("I16").Value = 1459.3846
$ws.Range("K16").Value = 1459.3846
$ws.Range("M16").Value = -1172.3846
$ws.Range("H31").Value = 2141.8438
$ws.Range("J31").Value = 2862.6365
$ws.Range("L31").Value = 2862.6365
$ws.Range("N31").Value = -3452.6365
$ws.Range("H34").Value = 2141.8438
$ws.Range("J34").Value = 2862.6365
$ws.Range("L34").Value = 2862.6365
$ws.Range("N34").Value = -3266.6365
$ws.Range("H58").Value = 1931.7
$ws.Range("I58").Value = 1799.8
$ws.Range("K58").Value = 1799.8
$ws.Range("M58").Value = -1596.8
$ws.Range("H105").Value = 1010
$ws.Range("I105").Value = 1012.73334
$ws.Range("J105").Value = 1004.1429
$ws.Range("K105").Value = 1012.73334
$ws.Range("L105").Value = 1004.1429
$ws.Range("M105").Value = 734.26666
$ws.Range("N105").Value = -4498.1429
$ws.Range("H113").Value = 1425.1428
$ws.Range("I113").Value = 1459.3846
$ws.Range("K113").Value = 1459.3846
$ws.Range("M113").Value = 710.6153999999999
$ws.Range("H136").Value = 1931.7
$ws.Range("I136").Value = 1799.8
$ws.Range("K136").Value = 5399.4
$ws.Range("M136").Value = -2849.4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 3313071.5
$ws.Range("I29").Value = 125043.25
$ws.Range("J29").Value = 8413917
$ws.Range("K29").Value = 375129.75
$ws.Range("L29").Value = 25241751
$ws.Range("M29").Value = -374852.75
$ws.Range("N29").Value = -25242305
$ws.Range("H38").Value = 239.2
$ws.Range("I38").Value = 236.75
$ws.Range("J38").Value = 249
$ws.Range("K38").Value = 710.25
$ws.Range("L38").Value = 747
$ws.Range("M38").Value = -363.25
$ws.Range("N38").Value = -1441
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 311.15384
$ws.Range("J2").Value = 385.14285
$ws.Range("L2").Value = 385.14285
$ws.Range("N2").Value = -611.14285
$ws.Range("H44").Value = 19250
$ws.Range("I44").Value = 10000
$ws.Range("J44").Value = 28500
$ws.Range("K44").Value = 10000
$ws.Range("L44").Value = 28500
$ws.Range("M44").Value = -9404
$ws.Range("N44").Value = -29692
$ws.Range("H80").Value = 3132.75
$ws.Range("I80").Value = 2600.5
$ws.Range("J80").Value = 4197.25
$ws.Range("K80").Value = 2600.5
$ws.Range("L80").Value = 4197.25
$ws.Range("M80").Value = -1602.5
$ws.Range("N80").Value = -6193.25
$ws.Range("H83").Value = 3132.75
$ws.Range("I83").Value = 2600.5
$ws.Range("J83").Value = 4197.25
$ws.Range("K83").Value = 13002.5
$ws.Range("L83").Value = 20986.25
$ws.Range("M83").Value = -8010.5
$ws.Range("N83").Value = -30970.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 85.27778000000001
$ws.Range("J2").Value = 85.27778000000001
$ws.Range("L2").Value = 85.27778000000001
$ws.Range("N2").Value = -309.27778
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H31").Value = 5437.95
$ws.Range("J31").Value = 8549.5
$ws.Range("L31").Value = 8549.5
$ws.Range("N31").Value = -9045.5
$ws.Range("H35").Value = 11999.5
$ws.Range("I35").Value = 3000
$ws.Range("J35").Value = 20999
$ws.Range("K35").Value = 3000
$ws.Range("L35").Value = 20999
$ws.Range("M35").Value = -2664
$ws.Range("N35").Value = -21671
$ws.Range("H40").Value = 5121.125
$ws.Range("I40").Value = 2908.7778
$ws.Range("K40").Value = 2908.7778
$ws.Range("M40").Value = -2772.7778
$ws.Range("H74").Value = 23108.5
$ws.Range("I74").Value = 21000
$ws.Range("K74").Value = 21000
$ws.Range("M74").Value = -20002
$ws.Range("H77").Value = 23108.5
$ws.Range("I77").Value = 21000
$ws.Range("K77").Value = 63000
$ws.Range("M77").Value = -58008
$ws.Range("H82").Value = 2997.2632
$ws.Range("I82").Value = 3397
$ws.Range("J82").Value = 1498.25
$ws.Range("K82").Value = 3397
$ws.Range("L82").Value = 1498.25
$ws.Range("M82").Value = -3036
$ws.Range("N82").Value = -2220.25
$ws.Range("H85").Value = 2997.2632
$ws.Range("I85").Value = 3397
$ws.Range("J85").Value = 1498.25
$ws.Range("K85").Value = 3397
$ws.Range("L85").Value = 1498.25
$ws.Range("M85").Value = -2149
$ws.Range("N85").Value = -3994.25
$ws.Range("H132").Value = 3146.4358
$ws.Range("I132").Value = 2546.577
$ws.Range("K132").Value = 7639.731000000001
$ws.Range("M132").Value = -5109.731000000001
$ws.Range("H136").Value = 3725.6538
$ws.Range("I136").Value = 2742
$ws.Range("K136").Value = 8226
$ws.Range("M136").Value = -5676
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 36247.5
$ws.Range("J44").Value = 36247.5
$ws.Range("L44").Value = 36247.5
$ws.Range("N44").Value = -37355.5
$ws.Range("H75").Value = 24008
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 24008
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H122").Value = 1704.4
$ws.Range("I122").Value = 1724
$ws.Range("J122").Value = 1634.1666
$ws.Range("K122").Value = 5172
$ws.Range("L122").Value = 4902.4998
$ws.Range("M122").Value = -2722
$ws.Range("N122").Value = -9802.4998
$ws.Range("H133").Value = 53427.832
$ws.Range("J133").Value = 53427.832
$ws.Range("L133").Value = 53427.832
$ws.Range("N133").Value = -63547.832
